# Add a new "MIDDLE HAND-LOWBALL" block in columns L:M (mirrors the existing
# BOTTOM/MIDDLE/TOP HAND blocks in A:B, E:F, I:J).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column widths for the two new columns.
# Target stored widths are L=22.85546875, M=12; the host quantizes
# ColumnWidth to 1/6-character steps, so these inputs land on the closest
# (M lands exactly on 12) achievable stored widths.
$ws.Columns.Item(12).ColumnWidth = 22
$ws.Columns.Item(13).ColumnWidth = 11.166666666666666

# Block title (row 1)
$ws.Range("L1").Value = "MIDDLE HAND-LOWBALL"

# Column headers (row 3)
$ws.Range("L3").Value = "HAND VALUE"
$ws.Range("M3").Value = "UNITS"

# Data rows (4-7)
$ws.Range("L4").Value = "Nine Low"
$ws.Range("M4").Value = 1

$ws.Range("L5").Value = "Eight Low"
$ws.Range("M5").Value = 2

$ws.Range("L6").Value = "Seven Low"
$ws.Range("M6").Value = 4

$ws.Range("L7").Value = "Wheel"
$ws.Range("M7").Value = 8

# Match the saved selection/view of the edited workbook
$ws.Range("M8").Select()
